$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new student record on row 8
$ws.Range("A8").Value = "21T2284"
$ws.Range("B8").Value = "ivan"
$ws.Range("C8").Value = "kamdem"
$ws.Range("D8").Value = "L3"
$ws.Range("E8").Value = "ivan.kamdem@facsciences-uy1.cm"
$ws.Range("F8").Value = "INFO"

# Add mailto hyperlink on the mail cell, mirroring the existing E2 hyperlink
$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:ivan.kamdem@facsciences-uy1.cm", "", "", "ivan.kamdem@facsciences-uy1.cm")

# Move the active selection to F8
$ws.Range("F8").Select()
